# Update gh-pages to output generated at 456a3b4
# Updates "想去人数" (want-to-go counts) and one event's name/cover image
# on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 4436
$ws1.Range("F4").Value = 337
$ws1.Range("F9").Value = 131
$ws1.Range("F10").Value = 314
$ws1.Range("F11").Value = 249
$ws1.Range("F12").Value = 2960

$ws1.Range("C14").Value = "江西·JMG（江西广电）第二届UP动漫游戏博览会"
$ws1.Range("F14").Value = 1536
$ws1.Range("I14").Value = "//i0.hdslb.com/bfs/openplatform/202408/oZpM885D1724642687206.png"

# ---- Sheet "全部类型" (all types) ----
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 4436
$ws4.Range("F4").Value = 337
$ws4.Range("F10").Value = 131
$ws4.Range("F11").Value = 314
$ws4.Range("F12").Value = 249
$ws4.Range("F13").Value = 2960

$ws4.Range("C15").Value = "江西·JMG（江西广电）第二届UP动漫游戏博览会"
$ws4.Range("F15").Value = 1536
$ws4.Range("I15").Value = "//i0.hdslb.com/bfs/openplatform/202408/oZpM885D1724642687206.png"
